$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (shared-string rich-text runs)
# "Volume 29   Number  38" -> "...40"  (run 4, chars 21-22)
# ---------------------------------------------------------------------------
$ws.Range("A8").Characters(21, 2).Text = "40"
$ws.Range("C9").Characters(27, 9).Text = "10/3/2022"
$ws.Range("C9").Characters(47, 9).Text = "10/9/2022"

# ---------------------------------------------------------------------------
# Type-swap cells (number <-> text) - use Copy() from a donor cell that already
# carries the desired style+type, then overwrite the value where it became numeric
# ---------------------------------------------------------------------------
$ws.Range("C14").Copy($ws.Range("F14"))
$ws.Range("C14").Copy($ws.Range("C15"))
$ws.Range("C14").Copy($ws.Range("G15"))
$ws.Range("E14").Copy($ws.Range("H15"))
$ws.Range("C14").Copy($ws.Range("D20"))
$ws.Range("E14").Copy($ws.Range("E20"))
$ws.Range("C14").Copy($ws.Range("D23"))
$ws.Range("E14").Copy($ws.Range("E23"))
$ws.Range("C14").Copy($ws.Range("D26"))
$ws.Range("E14").Copy($ws.Range("E26"))
$ws.Range("I14").Copy($ws.Range("C28"))
$ws.Range("C14").Copy($ws.Range("D28"))
$ws.Range("E14").Copy($ws.Range("E28"))
$ws.Range("I14").Copy($ws.Range("C29"))
$ws.Range("C14").Copy($ws.Range("D29"))
$ws.Range("E14").Copy($ws.Range("E29"))

# Overwrite the numeric value for cells that swapped text -> number
$ws.Range("C28").Value = 1
$ws.Range("C29").Value = 1

# ---------------------------------------------------------------------------
# Plain numeric value updates (style/type unchanged)
# ---------------------------------------------------------------------------
$ws.Range("N15").Value = -43.75
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 33.333333333333
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 18.181818181818
$ws.Range("I16").Value = 117
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 17
$ws.Range("L16").Value = 39.285714285714
$ws.Range("M16").Value = -13.333333333333
$ws.Range("N16").Value = -88.264794383149
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -44.444444444444
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 23
$ws.Range("H17").Value = -65.217391304347
$ws.Range("I17").Value = 131
$ws.Range("J17").Value = 132
$ws.Range("K17").Value = -0.757575757575
$ws.Range("L17").Value = 3.968253968253
$ws.Range("M17").Value = 48.863636363636
$ws.Range("N17").Value = -61.470588235294
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 50
$ws.Range("F18").Value = 11
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 149
$ws.Range("J18").Value = 104
$ws.Range("K18").Value = 43.26923076923
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 63.736263736263
$ws.Range("N18").Value = -76.273885350318
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = 77.777777777777
$ws.Range("F19").Value = 54
$ws.Range("G19").Value = 53
$ws.Range("H19").Value = 1.88679245283
$ws.Range("I19").Value = 498
$ws.Range("J19").Value = 397
$ws.Range("K19").Value = 25.44080604534
$ws.Range("L19").Value = 62.214983713355
$ws.Range("M19").Value = 30.708661417322
$ws.Range("N19").Value = -42.88990825688
$ws.Range("C20").Value = 1
$ws.Range("F20").Value = 7
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 16.666666666666
$ws.Range("I20").Value = 47
$ws.Range("J20").Value = 29
$ws.Range("K20").Value = 62.068965517241
$ws.Range("L20").Value = 30.555555555555
$ws.Range("M20").Value = 4.444444444444
$ws.Range("N20").Value = -91.710758377425
$ws.Range("C21").Value = 29
$ws.Range("D21").Value = 23
$ws.Range("E21").Value = 26.086956521739
$ws.Range("G21").Value = 104
$ws.Range("H21").Value = -7.692307692307
$ws.Range("I21").Value = 952
$ws.Range("J21").Value = 772
$ws.Range("K21").Value = 23.316062176165
$ws.Range("L21").Value = 34.463276836158
$ws.Range("M21").Value = 27.956989247311
$ws.Range("N21").Value = -72.244897959183
$ws.Range("C22").Value = 1
$ws.Range("D22").Value = 1
$ws.Range("G22").Value = 7
$ws.Range("H22").Value = -14.285714285714
$ws.Range("I22").Value = 50
$ws.Range("J22").Value = 39
$ws.Range("K22").Value = 28.205128205128
$ws.Range("L22").Value = 4.166666666666
$ws.Range("M22").Value = -10.714285714285
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 33.333333333333
$ws.Range("I23").Value = 29
$ws.Range("K23").Value = -21.621621621621
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 45
$ws.Range("C24").Value = 37
$ws.Range("D24").Value = 22
$ws.Range("E24").Value = 68.181818181818
$ws.Range("F24").Value = 196
$ws.Range("G24").Value = 107
$ws.Range("H24").Value = 83.177570093457
$ws.Range("I24").Value = 1486
$ws.Range("J24").Value = 1017
$ws.Range("K24").Value = 46.116027531956
$ws.Range("L24").Value = 102.176870748299
$ws.Range("M24").Value = 23.524522028262
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 19
$ws.Range("G25").Value = 20
$ws.Range("H25").Value = -5
$ws.Range("I25").Value = 281
$ws.Range("J25").Value = 198
$ws.Range("K25").Value = 41.919191919191
$ws.Range("L25").Value = 31.924882629108
$ws.Range("M25").Value = 4.074074074074
$ws.Range("C26").Value = 1
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 150
$ws.Range("I26").Value = 15
$ws.Range("J26").Value = 21
$ws.Range("K26").Value = -28.571428571428
$ws.Range("L26").Value = 200
$ws.Range("D27").Value = 3
$ws.Range("G27").Value = 10
$ws.Range("H27").Value = -70
$ws.Range("I27").Value = 51
$ws.Range("J27").Value = 52
$ws.Range("K27").Value = -1.923076923076
$ws.Range("L27").Value = 45.714285714285
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = -66.666666666666
$ws.Range("I28").Value = 6
$ws.Range("J28").Value = 9
$ws.Range("K28").Value = -33.333333333333
$ws.Range("L28").Value = -14.285714285714
$ws.Range("M28").Value = -25
$ws.Range("N28").Value = -79.310344827586
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = 3
$ws.Range("H29").Value = -66.666666666666
$ws.Range("I29").Value = 6
$ws.Range("J29").Value = 9
$ws.Range("K29").Value = -33.333333333333
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = 0
$ws.Range("N29").Value = -73.91304347826

